# Printing outputs for final paper
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename round labels: "Quarter" -> "Quarterfinals", "Semis" -> "Semifinals"
# These shared strings are used by cells A10:A13 ("Quarter") and A14:A15 ("Semis").
$ws.Range("A10:A13").Value = "Quarterfinals"
$ws.Range("A14:A15").Value = "Semifinals"
$ws.Range("A16").Value = "Finals"

# Fix the Quarterfinals row 10: the winner of Goffin D. vs Verdasco F. (row 3)
# was recorded incorrectly as Verdasco F./37; correct it to Goffin D./23.
$ws.Range("D10").Value = "Goffin D."
$ws.Range("E10").Value = 23

# Update the active selection left on the sheet.
$ws.Range("C21").Select()

$wb.Save()
